# Excel Sheet Downloading Functionality Done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Students Data"

# Header row
$ws.Range("A1").Value = "Fullname"
$ws.Range("B1").Value = "Rollno"
$ws.Range("C1").Value = "Mobileno"

# --- Apply the date-style cells first so the new custom number format
# (yyyy-mm-dd) claims cell-style index 1, matching the target workbook. ---
$ws.Range("O2").NumberFormat = "yyyy-mm-dd"
$ws.Range("O2").Value = (Get-Date -Year 2021 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("P2").NumberFormat = "yyyy-mm-dd"
$ws.Range("P2").Value = (Get-Date -Year 2021 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("O3").NumberFormat = "yyyy-mm-dd"
$ws.Range("O3").Value = (Get-Date -Year 2021 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("P3").NumberFormat = "yyyy-mm-dd"
$ws.Range("P3").Value = (Get-Date -Year 2021 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0)

# Row 2 data
$ws.Range("A2").Value = "Riya Ingale"
$ws.Range("B2").Value = "19102B0030"

# Mobile number must stay textual even though it looks numeric
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "8692931133"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "riya.ingale@vit.edu.in"
$ws.Range("E2").Value = "CMPN"
$ws.Range("F2").Value = "B"
$ws.Range("G2").Value = "TE"
$ws.Range("H2").Value = "Swabhav Techlabs"
$ws.Range("I2").Value = "Python Programming"
$ws.Range("J2").Value = "Self"
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = "Python"
$ws.Range("M2").Value = "Dhwani"

# N2 stays blank but the cell itself needs to exist
$ws.Range("N2").Font.Bold = $false

# Row 3 data
$ws.Range("A3").Value = "Samiksha"
$ws.Range("B3").Value = "19102B0021"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "9892805720"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "samiksha143018@gmail.com"
$ws.Range("E3").Value = "CMPN"
$ws.Range("F3").Value = "B"
$ws.Range("G3").Value = "TE"
$ws.Range("H3").Value = "Swabhav Techlabs"
$ws.Range("I3").Value = "Python Programming"
$ws.Range("J3").Value = "Self"
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = "Python"
$ws.Range("M3").Value = "Dhwani"

$ws.Range("N3").Font.Bold = $false
